$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Vista 3" requirement (row 4) as completed ("OK") in the Cumplido column (D)
$ws.Range("D4").Value = "OK"

# Move active selection to B7, matching the saved cursor position
$ws.Range("B7").Select()
